# ---------------------------------------------------------------------------
# Applies the "Add files via upload" commit:
#   * rename Sheet1 -> "latest"
#   * add a new sheet "peak" (after "latest")
#   * refresh several numeric values on "latest" (rows 2 & 3)
#   * populate "peak" with its header row + one data row of "last peak"
#     value/date pairs
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# -- 1. rename existing sheet, add the new one right after it --------------
$latest = $wb.Worksheets.Item(1)
$latest.Name = "latest"

$peak = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $latest)
$peak.Name = "peak"

# ---------------------------------------------------------------------------
# 2. update the changed values on "latest"
# ---------------------------------------------------------------------------

$row2 = @{
    "C2" = 17116;    "D2" = 35117;    "E2" = 45384;    "F2" = 73157;
    "G2" = 10874069; "H2" = 11953460; "I2" = 11668217; "J2" = 12779426;
    "K2" = 5420867;  "L2" = 7864483;  "M2" = 11537664; "N2" = 16243552;
    "O2" = 3262830336; "P2" = 3805640448;
    "AG2" = 1180652; "AH2" = 361040;
}
foreach ($addr in $row2.Keys) {
    $latest.Range($addr).Value = $row2[$addr]
}

# row 3 date moves from 44501 -> 44531
$latest.Range("A3").Value = 44531

$row3 = @{
    "C3" = 15263;    "D3" = 21401;    "E3" = 24276;    "F3" = 29505;
    "G3" = 11603034; "H3" = 13254094; "I3" = 13288900; "J3" = 14925248;
    "K3" = 4825758;  "L3" = 6766693;  "M3" = 6462811;  "N3" = 9365756;
    "O3" = 3492603648; "P3" = 4121775360;
    "AG3" = 983594;  "AH3" = 270364;
}
foreach ($addr in $row3.Keys) {
    $latest.Range($addr).Value = $row3[$addr]
}

# ---------------------------------------------------------------------------
# 3. populate the new "peak" sheet
# ---------------------------------------------------------------------------

$headers = @(
    "DayDeaMeSmA02S01_LPV", "DayDeaMeSmA02S01_LPD",
    "DayDeaUpSmA02S01_LPV", "DayDeaUpSmA02S01_LPD",
    "DayDeaMeSmA02S03_LPV", "DayDeaMeSmA02S03_LPD",
    "DayDeaUpSmA02S03_LPV", "DayDeaUpSmA02S03_LPD",
    "DayINFMeRaA02S01_LPV", "DayINFMeRaA02S01_LPD",
    "DayINFUpRaA02S01_LPV", "DayINFUpRaA02S01_LPD",
    "DayINFMeRaA02S03_LPV", "DayINFMeRaA02S03_LPD",
    "DayINFUpRaA02S03_LPV", "DayINFUpRaA02S03_LPD",
    "DayDeaMeSmA03S02_LPV", "DayDeaMeSmA03S02_LPD",
    "DayDeaUpRaA03S02_LPV", "DayDeaUpRaA03S02_LPD",
    "DayDeaMeRaA03S03_LPV", "DayDeaMeRaA03S03_LPD",
    "DayDeaUpRaA03S03_LPV", "DayDeaUpRaA03S03_LPD",
    "DayINFMeRaA03S02_LPV", "DayINFMeRaA03S02_LPD",
    "DayINFUpRaA03S02_LPV", "DayINFUpRaA03S02_LPD",
    "DayINFMeRaA03S03_LPV", "DayINFMeRaA03S03_LPD",
    "DayINFUpRaA03S03_LPV", "DayINFUpRaA03S03_LPD"
)

# value (odd / "LPV") columns hold a plain number, date ("LPD") columns hold
# an Excel date serial (re-using the same date style already used on
# "latest"!A2) - a few trailing LPD columns have no recorded date at all.
$values = @(
    26648, 44423,
    49705, 44468,
    45384, 44485,
    75746, 44477,
    8852569, 44422,
    18533956, 44442,
    14585255, 44461,
    23152494, 44455,
    15470, 44486,
    22432, $null,
    64595, $null,
    88491, $null,
    6546464, $null,
    9358429, $null,
    25558048, $null,
    33235220, $null
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $peak.Cells.Item(1, $col).Value = $headers[$i]
}

# grab the date-formatted style already present on "latest" so the LPD
# columns share the exact same number format (xf with numFmtId 14)
$latest.Range("A2").Copy()

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $isDateCol = ($i % 2) -eq 1
    $cell = $peak.Cells.Item(2, $col)
    if ($isDateCol) {
        $cell.PasteSpecial(-4122) # xlPasteFormats
        if ($null -ne $values[$i]) {
            $cell.Value = $values[$i]
        }
    } else {
        $cell.Value = $values[$i]
    }
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. leave the workbook selection on the first tab, as before
# ---------------------------------------------------------------------------
$latest.Activate()
